# refactor: update template layout
#
# The "Score Report" sheet had a title cell (A1, merged across A1:E1) that
# read "Survey:" and was left-aligned. Update it to read "Survey" (no
# trailing colon) and center it across the merged title range, then move
# the sheet's active selection onto that title range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 (merged A1:E1) title text: "Survey:" -> "Survey"
$ws.Range("A1").Value = "Survey"

# Center-align the title range (was left-aligned)
$ws.Range("A1:E1").HorizontalAlignment = -4108   # xlCenter

# Update the sheet's selection to the title range
$ws.Range("A1:E1").Select()
